# Update the evaluation report text values as per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 3. Số ngày làm việc theo quy định của pháp luật lao động trong tháng: 23 -> ...: 0
$ws.Range("A9").Value = "3. Số ngày làm việc theo quy định của pháp luật lao động trong tháng: 0"

# 7. Hành vi vi phạm: 0 -> 7. Hành vi vi phạm: (trailing space, value cleared)
$ws.Range("F12").Value = "7. Hành vi vi phạm: "

# 8. Hình thức kỷ luật: 0 -> 8. Hình thức kỷ luật: (trailing space, value cleared)
$ws.Range("I12").Value = "8. Hình thức kỷ luật: "
